$d = $word.ActiveDocument

# Locate the paragraph that ends in
#   "Organisationsname " + 76 underscores
# built out of three separate runs ("...underscores", "____", "_").
# The edit inserts an empty "_GoBack" bookmark right between the run
# that ends in "____" and the final lone "_" run - i.e. immediately
# after the 75th underscore following "Organisationsname ".
$underscores = "".PadLeft(75, '_')
$search = "Organisationsname " + $underscores

$rng = $d.Content
$found = $rng.Find.Execute($search, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    # Collapse the found range to its end point (right after the 75th
    # underscore / right after the "____" run) and drop the bookmark
    # there - this reproduces Word's automatic "_GoBack" bookmark that
    # marks the last edit position.
    $rng.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $rng)
    Write-Output "Inserted _GoBack bookmark"
} else {
    Write-Output "WARNING: anchor text not found; _GoBack bookmark NOT inserted"
}
